$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old content + formatting entirely, including the row that
# needs to disappear (old row 3).
$ws.Cells.Clear()

# Row 1 - header row
$ws.Range("A1").Value = "nama_survei"
$ws.Range("B1").Value = "lokasi_survei"
$ws.Range("C1").Value = "kode_desa"
$ws.Range("D1").Value = "kode_kecamatan"
$ws.Range("E1").Value = "kode_kabupaten"
$ws.Range("F1").Value = "kode_provinsi"
$ws.Range("G1").Value = "kro"
$ws.Range("H1").Value = "jadwal"
$ws.Range("I1").Value = "tim"

# V1 carries the text-formatted style but stays empty
$ws.Range("V1").NumberFormat = "@"

# Row 2 - data row
$ws.Range("A2").Value = "plan valencia"
$ws.Range("B2").Value = "london enggress"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "A"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "01-12-2029"
$ws.Range("I2").Value = 1

# Selection / active cell update
$ws.Range("N3").Select()
